$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: response changes from FALSE to "INSUFFICIENT INFO"; confidence is cleared
$ws.Range("B2").Value = "INSUFFICIENT INFO"
$ws.Range("C2").ClearContents()

# Row 3: confidence 0.8 -> 0.9
$ws.Range("C3").Value = 0.9

# Row 8: confidence 1 -> 0.9
$ws.Range("C8").Value = 0.9

# Row 9: response FALSE -> TRUE (kept as text, not boolean); confidence 1 -> 0.9
$ws.Range("B9").Value = "'TRUE"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 0.9
